# Auto-generated Excel COM-interop script
# Applies scheduled-runner market data sync updates to Sheets/Anima_Profits.xlsx
# (workbook has 8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 11
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()

# ALC!row 64
$ws.Range("H64").Value = 3033.2222
$ws.Range("I64").Value = 3060
$ws.Range("J64").Value = 2999.75
$ws.Range("K64").Value = 3060
$ws.Range("L64").Value = 2999.75
$ws.Range("M64").Value = -2812
$ws.Range("N64").Value = -3495.75

# ALC!row 67
$ws.Range("H67").Value = 3033.2222
$ws.Range("I67").Value = 3060
$ws.Range("J67").Value = 2999.75
$ws.Range("K67").Value = 3060
$ws.Range("L67").Value = 2999.75
$ws.Range("M67").Value = -2202
$ws.Range("N67").Value = -4715.75

# ALC!row 138
$ws.Range("H138").Value = 2709.9275
$ws.Range("I138").Value = 1843.5862
$ws.Range("J138").Value = 3338.025
$ws.Range("K138").Value = 5530.7586
$ws.Range("L138").Value = 10014.075
$ws.Range("M138").Value = -390.7586000000001
$ws.Range("N138").Value = -20294.075

# ALC!row 141
$ws.Range("H141").Value = 3352.282
$ws.Range("I141").Value = 1169.5161
$ws.Range("J141").Value = 11810.5
$ws.Range("K141").Value = 3508.5483
$ws.Range("L141").Value = 35431.5
$ws.Range("M141").Value = 1671.4517
$ws.Range("N141").Value = -45791.5

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 5
$ws.Range("H5").Value = 140
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 140
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 140
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -364

# ARM!row 32
$ws.Range("H32").Value = 864123.9399999999
$ws.Range("I32").Value = 944347.9
$ws.Range("J32").Value = 41828.5
$ws.Range("K32").Value = 944347.9
$ws.Range("L32").Value = 41828.5
$ws.Range("M32").Value = -944060.9
$ws.Range("N32").Value = -42402.5

# ARM!row 61
$ws.Range("H61").Value = 3053.3462
$ws.Range("I61").Value = 2534.158
$ws.Range("J61").Value = 4462.5713
$ws.Range("K61").Value = 2534.158
$ws.Range("L61").Value = 4462.5713
$ws.Range("M61").Value = -2322.158
$ws.Range("N61").Value = -4886.5713

# ARM!row 74
$ws.Range("H74").Value = 3000.7778
$ws.Range("I74").Value = 2750
$ws.Range("J74").Value = 3072.4285
$ws.Range("K74").Value = 2750
$ws.Range("L74").Value = 3072.4285
$ws.Range("M74").Value = -1876
$ws.Range("N74").Value = -4820.4285

# ARM!row 77
$ws.Range("H77").Value = 3000.7778
$ws.Range("I77").Value = 2750
$ws.Range("J77").Value = 3072.4285
$ws.Range("K77").Value = 13750
$ws.Range("L77").Value = 15362.1425
$ws.Range("M77").Value = -9382
$ws.Range("N77").Value = -24098.1425

# ARM!row 102
$ws.Range("H102").Value = 3684.6
$ws.Range("I102").Value = 3105.75
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 3105.75
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = -1483.75

# ARM!row 136
$ws.Range("H136").Value = 3053.3462
$ws.Range("I136").Value = 2534.158
$ws.Range("J136").Value = 4462.5713
$ws.Range("K136").Value = 7602.474
$ws.Range("L136").Value = 13387.7139
$ws.Range("M136").Value = -5052.474
$ws.Range("N136").Value = -18487.7139

$ws = $wb.Worksheets.Item("BSM")
# BSM!row 4
$ws.Range("H4").Value = 140
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 140
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 140
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -370

# BSM!row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()

# BSM!row 80
$ws.Range("H80").Value = 1189.1
$ws.Range("I80").Value = 2388.5557
$ws.Range("J80").Value = 207.72728
$ws.Range("K80").Value = 2388.5557
$ws.Range("L80").Value = 207.72728
$ws.Range("M80").Value = -1390.5557
$ws.Range("N80").Value = -2203.72728

# BSM!row 83
$ws.Range("H83").Value = 1189.1
$ws.Range("I83").Value = 2388.5557
$ws.Range("J83").Value = 207.72728
$ws.Range("K83").Value = 11942.7785
$ws.Range("L83").Value = 1038.6364
$ws.Range("M83").Value = -6950.7785
$ws.Range("N83").Value = -11022.6364

# BSM!row 99
$ws.Range("H99").Value = 1873.3334
$ws.Range("I99").Value = 1410
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 1410
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = 88

# BSM!row 107
$ws.Range("H107").Value = 91538
$ws.Range("I107").Value = 111499.78
$ws.Range("J107").Value = 1710
$ws.Range("K107").Value = 111499.78
$ws.Range("L107").Value = 1710
$ws.Range("M107").Value = -109579.78
$ws.Range("N107").Value = -5550

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 31
$ws.Range("H31").Value = 6806.1875
$ws.Range("I31").Value = 1291.8572
$ws.Range("J31").Value = 11095.111
$ws.Range("K31").Value = 1291.8572
$ws.Range("L31").Value = 11095.111
$ws.Range("M31").Value = -996.8571999999999
$ws.Range("N31").Value = -11685.111

# CRP!row 34
$ws.Range("H34").Value = 6806.1875
$ws.Range("I34").Value = 1291.8572
$ws.Range("J34").Value = 11095.111
$ws.Range("K34").Value = 1291.8572
$ws.Range("L34").Value = 11095.111
$ws.Range("M34").Value = -1089.8572
$ws.Range("N34").Value = -11499.111

# CRP!row 122
$ws.Range("H122").Value = 1864.1364
$ws.Range("I122").Value = 1202.2
$ws.Range("J122").Value = 2058.8235
$ws.Range("K122").Value = 3606.6
$ws.Range("L122").Value = 6176.470499999999
$ws.Range("M122").Value = -1156.6

# CRP!row 134
$ws.Range("H134").Value = 4934.7666
$ws.Range("I134").Value = 4736.077
$ws.Range("J134").Value = 6226.25
$ws.Range("K134").Value = 14208.231
$ws.Range("L134").Value = 18678.75
$ws.Range("M134").Value = -11673.231
$ws.Range("N134").Value = -23748.75

# CRP!row 141
$ws.Range("H141").Value = 216666.67
$ws.Range("I141").Value = 200000
$ws.Range("J141").Value = 220000
$ws.Range("K141").Value = 200000
$ws.Range("L141").Value = 220000
$ws.Range("M141").Value = -194820
$ws.Range("N141").Value = -230360

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 31
$ws.Range("H31").Value = 1929.8462
$ws.Range("I31").Value = 3000
$ws.Range("J31").Value = 1840.6666
$ws.Range("K31").Value = 9000
$ws.Range("L31").Value = 5521.9998
$ws.Range("M31").Value = -8712
$ws.Range("N31").Value = -6097.9998

# CUL!row 34
$ws.Range("H34").Value = 10204579
$ws.Range("I34").Value = 161
$ws.Range("J34").Value = 10870085
$ws.Range("K34").Value = 483
$ws.Range("L34").Value = 32610255
$ws.Range("M34").Value = -399
$ws.Range("N34").Value = -32610423

# CUL!row 35
$ws.Range("H35").Value = 4776.5
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 4776.5
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 14329.5
$ws.Range("N35").Value = -14905.5

# CUL!row 39
$ws.Range("H39").Value = 2361.5386
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2361.5386
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 7084.6158
$ws.Range("N39").Value = -7672.6158

# CUL!row 55
$ws.Range("H55").Value = 1658.3334
$ws.Range("I55").Value = 200
$ws.Range("J55").Value = 1790.909
$ws.Range("K55").Value = 600
$ws.Range("L55").Value = 5372.727000000001
$ws.Range("M55").Value = -423
$ws.Range("N55").Value = -5726.727000000001

# CUL!row 68
$ws.Range("H68").Value = 862.75
$ws.Range("I68").Value = 934
$ws.Range("J68").Value = 820
$ws.Range("K68").Value = 2802
$ws.Range("L68").Value = 2460
$ws.Range("M68").Value = -1991
$ws.Range("N68").Value = -4082

# CUL!row 71
$ws.Range("H71").Value = 862.75
$ws.Range("I71").Value = 934
$ws.Range("J71").Value = 820
$ws.Range("K71").Value = 8406
$ws.Range("L71").Value = 7380
$ws.Range("M71").Value = -4350
$ws.Range("N71").Value = -15492

# CUL!row 92
$ws.Range("H92").Value = 430.2
$ws.Range("I92").Value = 449.5
$ws.Range("J92").Value = 425.375
$ws.Range("K92").Value = 1348.5
$ws.Range("L92").Value = 1276.125
$ws.Range("M92").Value = -100.5
$ws.Range("N92").Value = -3772.125

# CUL!row 131
$ws.Range("H131").Value = 1125.5588
$ws.Range("I131").Value = 732.25
$ws.Range("J131").Value = 1178
$ws.Range("K131").Value = 2196.75
$ws.Range("L131").Value = 3534
$ws.Range("M131").Value = 2843.25
$ws.Range("N131").Value = -13614

# CUL!row 137
$ws.Range("H137").Value = 9269496
$ws.Range("I137").Value = 55589224
$ws.Range("J137").Value = 5550.467
$ws.Range("K137").Value = 166767672
$ws.Range("L137").Value = 16651.401
$ws.Range("M137").Value = -166762572
$ws.Range("N137").Value = -26851.401

# CUL!row 139
$ws.Range("H139").Value = 3961.52
$ws.Range("I139").Value = 1143.3334
$ws.Range("J139").Value = 6562.923
$ws.Range("K139").Value = 3430.0002
$ws.Range("L139").Value = 19688.769
$ws.Range("M139").Value = 1709.9998

# CUL!row 140
$ws.Range("H140").Value = 2022.3529
$ws.Range("I140").Value = 1709.091
$ws.Range("J140").Value = 2596.6667
$ws.Range("K140").Value = 5127.272999999999
$ws.Range("L140").Value = 7790.000100000001
$ws.Range("M140").Value = 52.72700000000077
$ws.Range("N140").Value = -18150.0001

$ws = $wb.Worksheets.Item("GSM")
# GSM!row 102
$ws.Range("H102").Value = 1629.1
$ws.Range("I102").Value = 1517.8667
$ws.Range("J102").Value = 1962.8
$ws.Range("K102").Value = 1517.8667
$ws.Range("L102").Value = 1962.8
$ws.Range("M102").Value = 104.1333

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 7
$ws.Range("H7").Value = 4600.385
$ws.Range("I7").Value = 3128.5715
$ws.Range("J7").Value = 6317.5
$ws.Range("K7").Value = 3128.5715
$ws.Range("L7").Value = 6317.5
$ws.Range("M7").Value = -3016.5715
$ws.Range("N7").Value = -6541.5

# LTW!row 40
$ws.Range("H40").Value = 203400.8
$ws.Range("I40").Value = 336334.66
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 336334.66
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -336198.66

# LTW!row 55
$ws.Range("H55").Value = 647.3333
$ws.Range("I55").Value = 531.3333
$ws.Range("J55").Value = 763.3333
$ws.Range("K55").Value = 531.3333
$ws.Range("L55").Value = 763.3333
$ws.Range("M55").Value = -358.3333
$ws.Range("N55").Value = -1109.3333

# LTW!row 126
$ws.Range("H126").Value = 4600.385
$ws.Range("I126").Value = 3128.5715
$ws.Range("J126").Value = 6317.5
$ws.Range("K126").Value = 9385.7145
$ws.Range("L126").Value = 18952.5
$ws.Range("M126").Value = -6915.7145
$ws.Range("N126").Value = -23892.5

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 122
$ws.Range("H122").Value = 2340.7407
$ws.Range("I122").Value = 2267.7273
$ws.Range("J122").Value = 2662
$ws.Range("K122").Value = 6803.1819
$ws.Range("L122").Value = 7986
$ws.Range("M122").Value = -4353.1819

# WVR!row 136
$ws.Range("H136").Value = 4552.8335
$ws.Range("I136").Value = 3904.25
$ws.Range("J136").Value = 5850
$ws.Range("K136").Value = 11712.75
$ws.Range("L136").Value = 17550
$ws.Range("M136").Value = -9162.75
$ws.Range("N136").Value = -22650
